$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells for columns G and H on row 3
$ws.Range("G3").Value = "RPM à vide "
$ws.Range("H3").Value = "VBAT à VIDE"

# Row 9
$ws.Range("A9").Value = 7.5
$ws.Range("B9").Value = 200
$ws.Range("C9").Value = 4.6500000000000004
$ws.Range("D9").Formula = "=Distance/C9"
$ws.Range("E9").Formula = "=D9/(Diam_roue*PI())*1000*60"
$ws.Range("F9").Formula = "=D9*3.6"
$ws.Range("G9").Value = 216
$ws.Range("H9").Value = 8.07

# Row 10
$ws.Range("A10").Value = 7.53
$ws.Range("B10").Value = 200
$ws.Range("C10").Value = 4.51
$ws.Range("D10").Formula = "=Distance/C10"
$ws.Range("E10").Formula = "=D10/(Diam_roue*PI())*1000*60"
$ws.Range("F10").Formula = "=D10*3.6"
$ws.Range("G10").Value = 217
$ws.Range("H10").Value = 7.87

# Row 11 - only A11 populated (used range keeps growing for the detection test)
$ws.Range("A11").Value = 7.2

# Copy the number formats (with fill) from the row above onto the new formula cells
$ws.Range("D4:F4").Copy()
$ws.Range("D9:F9").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$ws.Range("D4:F4").Copy()
$ws.Range("D10:F10").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)

$excel.CutCopyMode = $false

# Move the selection/view back to A12 (matches author re-selecting after editing row 11)
$ws.Range("A12").Select()
